$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '97.884.69'
$ws.Cells.Item(2, 5).Value = '  +3.05%  '
$ws.Cells.Item(3, 4).Value = '3.600.75'
$ws.Cells.Item(3, 5).Value = '  +1.67%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '242.38'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +3.17%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '657.35'
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +15.97%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.419'
$c.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +5.82%  '
$ws.Cells.Item(9, 5).Value = '  +6.86%  '
$ws.Cells.Item(10, 5).Value = '  -0.05%  '
$ws.Cells.Item(11, 4).Value = '3.597.14'
$ws.Cells.Item(11, 5).Value = '  +1.60%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '44.46'
$c.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +6.07%  '
$ws.Cells.Item(13, 5).Value = '  +1.68%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '6.46'
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +0.72%  '
$ws.Cells.Item(15, 4).Value = '97.819.38'
$ws.Cells.Item(15, 5).Value = '  +3.11%  '
$ws.Cells.Item(16, 4).Value = '4.269.38'
$ws.Cells.Item(16, 5).Value = '  +1.65%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000259'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +3.22%  '
$ws.Cells.Item(18, 5).Value = '  +9.68%  '
$ws.Cells.Item(19, 4).Value = '3.588.09'
$ws.Cells.Item(19, 5).Value = '  +1.22%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '12.73'
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +1.49%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '18.02'
$c.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +2.15%  '
$ws.Cells.Item(22, 5).Value = '  +10.40%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '3.49'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +0.86%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '515.47'
$c.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +2.85%  '
$ws.Cells.Item(25, 5).Value = '  +6.85%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '6.83'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +3.07%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '101.02'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +10.34%  '
$ws.Cells.Item(28, 5).Value = '  +5.43%  '
$ws.Cells.Item(29, 4).Value = '3.794.17'
$ws.Cells.Item(29, 5).Value = '  +1.73%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '0.158'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +13.82%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -0.05%  '
$ws.Cells.Item(32, 5).Value = '  +4.35%  '
$ws.Cells.Item(33, 5).Value = '  -0.26%  '
$ws.Cells.Item(34, 5).Value = '  +4.82%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -0.59%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '31.72'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +0.73%  '
$ws.Cells.Item(37, 2).Value = 'Bittensor'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '620.42'
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +9.41%  '
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '8.82'
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +6.49%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.570'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +3.12%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '1.63'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +4.02%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +14.40%  '
$ws.Cells.Item(42, 5).Value = '  +3.39%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.923'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +3.36%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '5.97'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +7.41%  '
$ws.Cells.Item(46, 5).Value = '  +8.56%  '
$ws.Cells.Item(47, 5).Value = '  +1.89%  '
$ws.Cells.Item(48, 5).Value = '  +1.29%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '33.34'
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +0.09%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '8.49'
$c.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +6.74%  '
$ws.Cells.Item(51, 2).Value = 'MantraDAO'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '3.57'
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +0.50%  '
